# db uploader 1차 완료
# Change the DATA_TYPE for the "id" row from "integer" to "SERIAL",
# and update the active selection/view to C3 (scrolling back so A1 is
# the top-left cell again, matching the saved sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "SERIAL"

[void]$ws.Range("A1").Select()
[void]$ws.Range("C3").Select()
